$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.040319673023562
$ws.Cells.Item(2, 4).Value = 1.050398476648657
$ws.Cells.Item(2, 5).Value = 1.038734453878905
$ws.Cells.Item(2, 6).Value = 1.061154357488823
$ws.Cells.Item(2, 9).Value = 1.046005074714861
$ws.Cells.Item(2, 10).Value = 1.045406919576044
$ws.Cells.Item(2, 11).Value = 1.053152864107349
$ws.Cells.Item(2, 12).Value = 1.041521643239324
$ws.Cells.Item(2, 13).Value = 1.063879189130017
$ws.Cells.Item(2, 14).Value = 1.046891517311692

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.041212101358133
$ws.Cells.Item(3, 4).Value = 1.051133711407257
$ws.Cells.Item(3, 5).Value = 1.039492228859288
$ws.Cells.Item(3, 6).Value = 1.062054522704025
$ws.Cells.Item(3, 9).Value = 1.046279186495366
$ws.Cells.Item(3, 10).Value = 1.045945201768942
$ws.Cells.Item(3, 11).Value = 1.053700736351665
$ws.Cells.Item(3, 12).Value = 1.042089630534628
$ws.Cells.Item(3, 13).Value = 1.064593704008221
$ws.Cells.Item(3, 14).Value = 1.047430563927046

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.04179006716735
$ws.Cells.Item(4, 4).Value = 1.051609868137562
$ws.Cells.Item(4, 5).Value = 1.03998335960324
$ws.Cells.Item(4, 6).Value = 1.062637770506372
$ws.Cells.Item(4, 9).Value = 1.046455524106837
$ws.Cells.Item(4, 10).Value = 1.046293362160718
$ws.Cells.Item(4, 11).Value = 1.054054987625185
$ws.Cells.Item(4, 12).Value = 1.042457292355482
$ws.Cells.Item(4, 13).Value = 1.065056186371493
$ws.Cells.Item(4, 14).Value = 1.047779218746512

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.042033163598755
$ws.Cells.Item(5, 4).Value = 1.051810141202698
$ws.Cells.Item(5, 5).Value = 1.040190020843837
$ws.Cells.Item(5, 6).Value = 1.062883153162266
$ws.Cells.Item(5, 9).Value = 1.046529408887993
$ws.Cells.Item(5, 10).Value = 1.046439693174671
$ws.Cells.Item(5, 11).Value = 1.054203851754462
$ws.Cells.Item(5, 12).Value = 1.042611888768136
$ws.Cells.Item(5, 13).Value = 1.065250646906174
$ws.Cells.Item(5, 14).Value = 1.04792575756729

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.042073987535461
$ws.Cells.Item(6, 4).Value = 1.051843773561738
$ws.Cells.Item(6, 5).Value = 1.040224731264401
$ws.Cells.Item(6, 6).Value = 1.062924364807607
$ws.Cells.Item(6, 9).Value = 1.046541799930544
$ws.Cells.Item(6, 10).Value = 1.04646426067545
$ws.Cells.Item(6, 11).Value = 1.05422884296218
$ws.Cells.Item(6, 12).Value = 1.042637847991813
$ws.Cells.Item(6, 13).Value = 1.065283299577237
$ws.Cells.Item(6, 14).Value = 1.047950359956738

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.041793314965409
$ws.Cells.Item(7, 4).Value = 1.051612543817733
$ws.Cells.Item(7, 5).Value = 1.039986120276734
$ws.Cells.Item(7, 6).Value = 1.062641048595443
$ws.Cells.Item(7, 9).Value = 1.046456512332498
$ws.Cells.Item(7, 10).Value = 1.046295317584203
$ws.Cells.Item(7, 11).Value = 1.054056977004398
$ws.Cells.Item(7, 12).Value = 1.042459357958493
$ws.Cells.Item(7, 13).Value = 1.065058784634955
$ws.Cells.Item(7, 14).Value = 1.047781176946923

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.040621168015003
$ws.Cells.Item(8, 4).Value = 1.050646866717535
$ws.Cells.Item(8, 5).Value = 1.038990380965924
$ws.Cells.Item(8, 6).Value = 1.061458409908052
$ws.Cells.Item(8, 9).Value = 1.046097925126047
$ws.Cells.Item(8, 10).Value = 1.045588863756852
$ws.Cells.Item(8, 11).Value = 1.053338072958943
$ws.Cells.Item(8, 12).Value = 1.041713568317314
$ws.Cells.Item(8, 13).Value = 1.06412063223521
$ws.Cells.Item(8, 14).Value = 1.047073719874106

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.038559620238294
$ws.Cells.Item(9, 4).Value = 1.048948433674453
$ws.Cells.Item(9, 5).Value = 1.037241954644436
$ws.Cells.Item(9, 6).Value = 1.059380500622043
$ws.Cells.Item(9, 9).Value = 1.045458183869917
$ws.Cells.Item(9, 10).Value = 1.044342946097775
$ws.Cells.Item(9, 11).Value = 1.052069342602357
$ws.Cells.Item(9, 12).Value = 1.040400486321973
$ws.Cells.Item(9, 13).Value = 1.06246864822727
$ws.Cells.Item(9, 14).Value = 1.045826032869044

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.037187970294092
$ws.Cells.Item(10, 4).Value = 1.047818396933361
$ws.Cells.Item(10, 5).Value = 1.036080594237994
$ws.Cells.Item(10, 6).Value = 1.057999392598542
$ws.Cells.Item(10, 9).Value = 1.045026444389208
$ws.Cells.Item(10, 10).Value = 1.043511689505595
$ws.Cells.Item(10, 11).Value = 1.051222294470174
$ws.Cells.Item(10, 12).Value = 1.039525904883385
$ws.Cells.Item(10, 13).Value = 1.061368188424291
$ws.Cells.Item(10, 14).Value = 1.044993595797157

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.036594690591677
$ws.Cells.Item(11, 4).Value = 1.047329631861581
$ws.Cells.Item(11, 5).Value = 1.0355787416243
$ws.Cells.Item(11, 6).Value = 1.057402365587989
$ws.Cells.Item(11, 9).Value = 1.044838261169863
$ws.Cells.Item(11, 10).Value = 1.043151606863876
$ws.Cells.Item(11, 11).Value = 1.050855236555146
$ws.Cells.Item(11, 12).Value = 1.039147407557406
$ws.Cells.Item(11, 13).Value = 1.060891899134582
$ws.Cells.Item(11, 14).Value = 1.044633001796785

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.036374419560971
$ws.Cells.Item(12, 4).Value = 1.047148166660386
$ws.Cells.Item(12, 5).Value = 1.035392486649263
$ws.Cells.Item(12, 6).Value = 1.057180755304115
$ws.Cells.Item(12, 9).Value = 1.044768176357763
$ws.Cells.Item(12, 10).Value = 1.043017835720082
$ws.Cells.Item(12, 11).Value = 1.050718853910768
$ws.Cells.Item(12, 12).Value = 1.039006848355627
$ws.Cells.Item(12, 13).Value = 1.060715018070277
$ws.Cells.Item(12, 14).Value = 1.044499040682622

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.036421663923693
$ws.Cells.Item(13, 4).Value = 1.047187087740373
$ws.Cells.Item(13, 5).Value = 1.035432431914896
$ws.Cells.Item(13, 6).Value = 1.057228284552197
$ws.Cells.Item(13, 9).Value = 1.044783218156865
$ws.Cells.Item(13, 10).Value = 1.043046530991224
$ws.Cells.Item(13, 11).Value = 1.050748110287456
$ws.Cells.Item(13, 12).Value = 1.039036997334577
$ws.Cells.Item(13, 13).Value = 1.060752958099073
$ws.Cells.Item(13, 14).Value = 1.044527776704341

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.036576480883266
$ws.Cells.Item(14, 4).Value = 1.04731463017086
$ws.Cells.Item(14, 5).Value = 1.035563342550144
$ws.Cells.Item(14, 6).Value = 1.057384044100061
$ws.Cells.Item(14, 9).Value = 1.044832471711185
$ws.Cells.Item(14, 10).Value = 1.043140549709705
$ws.Cells.Item(14, 11).Value = 1.050843963952787
$ws.Cells.Item(14, 12).Value = 1.039135788239357
$ws.Cells.Item(14, 13).Value = 1.060877277389406
$ws.Cells.Item(14, 14).Value = 1.044621928940187

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.036671881939243
$ws.Cells.Item(15, 4).Value = 1.047393224436308
$ws.Cells.Item(15, 5).Value = 1.035644021556146
$ws.Cells.Item(15, 6).Value = 1.057480032905155
$ws.Cells.Item(15, 9).Value = 1.044862793933251
$ws.Cells.Item(15, 10).Value = 1.043198475071736
$ws.Cells.Item(15, 11).Value = 1.050903017162161
$ws.Cells.Item(15, 12).Value = 1.039196660793251
$ws.Cells.Item(15, 13).Value = 1.060953879151625
$ws.Cells.Item(15, 14).Value = 1.044679936562876

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.037227358194416
$ws.Cells.Item(16, 4).Value = 1.047850846344321
$ws.Cells.Item(16, 5).Value = 1.036113922253952
$ws.Cells.Item(16, 6).Value = 1.058039036560751
$ws.Cells.Item(16, 9).Value = 1.045038907482351
$ws.Cells.Item(16, 10).Value = 1.043535584082328
$ws.Cells.Item(16, 11).Value = 1.051246649083614
$ws.Cells.Item(16, 12).Value = 1.039551028865217
$ws.Cells.Item(16, 13).Value = 1.061399802866329
$ws.Cells.Item(16, 14).Value = 1.045017524306929

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.03757596968494
$ws.Cells.Item(17, 4).Value = 1.048138048353762
$ws.Cells.Item(17, 5).Value = 1.036408953838165
$ws.Cells.Item(17, 6).Value = 1.058389954152218
$ws.Cells.Item(17, 9).Value = 1.045149048115449
$ws.Cells.Item(17, 10).Value = 1.043747006117105
$ws.Cells.Item(17, 11).Value = 1.051462126025878
$ws.Cells.Item(17, 12).Value = 1.039773369563007
$ws.Cells.Item(17, 13).Value = 1.061679578115081
$ws.Cells.Item(17, 14).Value = 1.045229246585246

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.037779371865454
$ws.Cells.Item(18, 4).Value = 1.048305621226405
$ws.Cells.Item(18, 5).Value = 1.036581139475827
$ws.Cells.Item(18, 6).Value = 1.0585947350049
$ws.Cells.Item(18, 9).Value = 1.045213171867704
$ws.Cells.Item(18, 10).Value = 1.0438703111385
$ws.Cells.Item(18, 11).Value = 1.051587782941015
$ws.Cells.Item(18, 12).Value = 1.039903076561463
$ws.Cells.Item(18, 13).Value = 1.061842787135924
$ws.Cells.Item(18, 14).Value = 1.045352726713915

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.037848737400071
$ws.Cells.Item(19, 4).Value = 1.048362768154569
$ws.Cells.Item(19, 5).Value = 1.036639867023079
$ws.Cells.Item(19, 6).Value = 1.05866457631272
$ws.Cells.Item(19, 9).Value = 1.045235016134326
$ws.Cells.Item(19, 10).Value = 1.043912352578296
$ws.Cells.Item(19, 11).Value = 1.051630624063536
$ws.Cells.Item(19, 12).Value = 1.039947306556027
$ws.Cells.Item(19, 13).Value = 1.061898440689392
$ws.Cells.Item(19, 14).Value = 1.045394827857377

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.037538560446246
$ws.Cells.Item(20, 4).Value = 1.048107228823808
$ws.Cells.Item(20, 5).Value = 1.036377289520348
$ws.Cells.Item(20, 6).Value = 1.058352294016559
$ws.Cells.Item(20, 9).Value = 1.045137243414636
$ws.Cells.Item(20, 10).Value = 1.043724323960993
$ws.Cells.Item(20, 11).Value = 1.051439010188114
$ws.Cells.Item(20, 12).Value = 1.039749512491977
$ws.Cells.Item(20, 13).Value = 1.061649558711947
$ws.Cells.Item(20, 14).Value = 1.045206532217871

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.036530888391604
$ws.Cells.Item(21, 4).Value = 1.047277069775386
$ws.Cells.Item(21, 5).Value = 1.035524788327671
$ws.Cells.Item(21, 6).Value = 1.057338172582939
$ws.Cells.Item(21, 9).Value = 1.044817972873557
$ws.Cells.Item(21, 10).Value = 1.043112864095473
$ws.Cells.Item(21, 11).Value = 1.050815738559776
$ws.Cells.Item(21, 12).Value = 1.039106695899587
$ws.Cells.Item(21, 13).Value = 1.060840667507804
$ws.Cells.Item(21, 14).Value = 1.044594204009206

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.03589790107844
$ws.Cells.Item(22, 4).Value = 1.046755602061771
$ws.Cells.Item(22, 5).Value = 1.034989687042863
$ws.Cells.Item(22, 6).Value = 1.056701435165168
$ws.Cells.Item(22, 9).Value = 1.044616163727308
$ws.Cells.Item(22, 10).Value = 1.042728297409636
$ws.Cells.Item(22, 11).Value = 1.050423626354568
$ws.Cells.Item(22, 12).Value = 1.03870271476114
$ws.Cells.Item(22, 13).Value = 1.060332282808474
$ws.Cells.Item(22, 14).Value = 1.044209091194565

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.036233404605212
$ws.Cells.Item(23, 4).Value = 1.047031995362303
$ws.Cells.Item(23, 5).Value = 1.035273268519846
$ws.Cells.Item(23, 6).Value = 1.05703889767339
$ws.Cells.Item(23, 9).Value = 1.044723247870887
$ws.Cells.Item(23, 10).Value = 1.042932174274199
$ws.Cells.Item(23, 11).Value = 1.050631514496352
$ws.Cells.Item(23, 12).Value = 1.038916855050639
$ws.Cells.Item(23, 13).Value = 1.060601767995467
$ws.Cells.Item(23, 14).Value = 1.044413257587661

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.037555463876421
$ws.Cells.Item(24, 4).Value = 1.048121154679817
$ws.Cells.Item(24, 5).Value = 1.036391596957796
$ws.Cells.Item(24, 6).Value = 1.05836931071263
$ws.Cells.Item(24, 9).Value = 1.045142577819848
$ws.Cells.Item(24, 10).Value = 1.043734573093982
$ws.Cells.Item(24, 11).Value = 1.05144945532417
$ws.Cells.Item(24, 12).Value = 1.039760292415553
$ws.Cells.Item(24, 13).Value = 1.061663123123773
$ws.Cells.Item(24, 14).Value = 1.045216795905804

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.03909210677885
$ws.Cells.Item(25, 4).Value = 1.049387129099385
$ws.Cells.Item(25, 5).Value = 1.037693221209451
$ws.Cells.Item(25, 6).Value = 1.059916962959031
$ws.Cells.Item(25, 9).Value = 1.045624499962655
$ws.Cells.Item(25, 10).Value = 1.044665164094534
$ws.Cells.Item(25, 11).Value = 1.052397561323575
$ws.Cells.Item(25, 12).Value = 1.040739812499599
$ws.Cells.Item(25, 13).Value = 1.06289557966364
$ws.Cells.Item(25, 14).Value = 1.046148708452317
